$d = $word.ActiveDocument

# 1) "Programa" section - Portuguese paragraph
$d.Content.Find.Execute(
    "- Fundamentos de análise titulométrica e cálculos em análise titulométrica.- Titulometria de neutralização: fundamentos, indicadores de titulação, curvas de titulação ácido base.- Titulometria complexométrica: fundamentos, complexometria com EDTA.- Titulação de oxiredução: fundamentos e principais indicadores- Titulometria de precipitação: fundamentos, indicadores, argentimetria.A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Fundamentos de análise titulométrica e cálculos em análise titulométrica.^l- Titulometria de neutralização: fundamentos, indicadores de titulação, curvas de titulação ácido base.^l- Titulometria complexométrica: fundamentos, complexometria com EDTA.^l- Titulação de oxiredução: fundamentos e principais indicadores^l- Titulometria de precipitação: fundamentos, indicadores, argentimetria.^lA disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina",
    2)

# 2) "Programa" section - English paragraph
$d.Content.Find.Execute(
    "- Fundamentals of titulometric analysis and calculations in titulometric analysis.- Neutralization titrometry: fundamentals, titration indicators, acid base titration curves.- Complexometric titrometry: fundamentals, complexometry with EDTA.- Titration of oxireduction: fundamentals and main indicators- Precipitation titrometry: fundamentals, indicators, argentimetry.The discipline may have didactic trips to complement the content of the discipline",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- Fundamentals of titulometric analysis and calculations in titulometric analysis.^l- Neutralization titrometry: fundamentals, titration indicators, acid base titration curves.^l- Complexometric titrometry: fundamentals, complexometry with EDTA.^l- Titration of oxireduction: fundamentals and main indicators^l- Precipitation titrometry: fundamentals, indicators, argentimetry.^lThe discipline may have didactic trips to complement the content of the discipline",
    2)

# 3) "Bibliografia" section
$d.Content.Find.Execute(
    "Harris, D.C. EXPLORANDO A QUÍMICA ANALÍTICA, 4ª edição, LTC, Rio de Janeiro – RJ, 2011Skoog, D.A., Holler, F.J. e Nieman, T.A., PRINCÍPIOS DE ANÁLISE INSTRUMENTAL, 5ª ed., Bookman, Porto Alegre, 2002.Mendham,J., Denney, R.C., Barnes, J.D. e Thomas, M., Vogel: ANÁLISE QUÍMICA QUANTITATIVA, 6ª ed., Livros Técnicos e Científicos, Rio de Janeiro -RJ, 2002.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Harris, D.C. EXPLORANDO A QUÍMICA ANALÍTICA, 4ª edição, LTC, Rio de Janeiro – RJ, 2011^lSkoog, D.A., Holler, F.J. e Nieman, T.A., PRINCÍPIOS DE ANÁLISE INSTRUMENTAL, 5ª ed., Bookman, Porto Alegre, 2002.^lMendham,J., Denney, R.C., Barnes, J.D. e Thomas, M., Vogel: ANÁLISE QUÍMICA QUANTITATIVA, 6ª ed., Livros Técnicos e Científicos, Rio de Janeiro -RJ, 2002.",
    2)

Write-Host "Done"
